$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualiza os dados do último trimestre (linha 21: 2025Q3)
$ws.Range("C21").Value = 295
$ws.Range("D21").Value = 250
$ws.Range("E21").Value = 45
$ws.Range("F21").Value = 71.63323782234957
